$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-PackageXml($bodyInnerXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyInnerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- Edit 1: insert a new, empty bold paragraph right after the
#     "Task 1: Proposed Solutions for Leveraging Azure Cloud" heading. ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Task 1: Proposed Solutions for Leveraging Azure Cloud", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$rng1.Collapse(0)
$rng1.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$newParaRange = $newPara.Range
$emptyParaXml = '<w:body><w:p ' + $wns + '><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p></w:body>'
$newParaRange.InsertXML((New-PackageXml $emptyParaXml))

# --- Edit 2: split the "Use Azure Cost Management..." run in two, moving
#     the lastRenderedPageBreak onto the start of the second half. ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Use Azure Cost Management", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$rng2.Expand(4)
$splitParaInner = '<w:r ' + $wns + '><w:t xml:space="preserve">Use Azure Cost Management + Billing to monitor and manage costs appropriately. This tool will be used to give insight into spending patterns and </w:t></w:r><w:r ' + $wns + '><w:lastRenderedPageBreak/><w:t>to cost the bill accurately according to the various departments or projects that are being processed.</w:t></w:r>'
$rng2.InsertXML((New-PackageXml ('<w:body><w:p>' + $splitParaInner + '</w:p></w:body>')))

# --- Edit 3: remove the (now redundant) lastRenderedPageBreak from the
#     start of the following "Create budgets..." run. ---
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("Create budgets for every resource group", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$rng3.Expand(4)
$budgetParaInner = '<w:r ' + $wns + '><w:t>Create budgets for every resource group and allow alerting when consumption is close to budget limits.</w:t></w:r>'
$rng3.InsertXML((New-PackageXml ('<w:body><w:p>' + $budgetParaInner + '</w:p></w:body>')))
